{"js": "// Add an \"Author: Craig Jones\" paragraph right after the \"Line.java\"\n// title paragraph (and before the existing blank paragraph that\n// follows it), matching the new-paragraph formatting (NoSpacing style,\n// Times New Roman ascii/hAnsi/cs fonts) already used by that blank\n// paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the title paragraph whose text is exactly \"Line.java\".\nlet titleIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Line.java\") {\n    titleIndex = i;\n    break;\n  }\n}\n\nlet anchorParagraph;\nlet insertLocation;\nlet needsExplicitFormatting = false;\n\nif (titleIndex !== -1 && titleIndex + 1 < paragraphs.items.length) {\n  // Insert right before the paragraph that currently follows the title\n  // (normally the blank separator paragraph) so the new paragraph picks\n  // up that paragraph's formatting (No Spacing / Times New Roman, no\n  // bold/large title sizing) exactly like Word does when typing a new\n  // line in front of existing, identically-styled text.\n  anchorParagraph = paragraphs.items[titleIndex + 1];\n  insertLocation = \"Before\";\n} else if (titleIndex !== -1) {\n  // Title paragraph is the very last paragraph in the document - fall\n  // back to inserting straight after it, then fix up formatting\n  // explicitly since it would otherwise inherit the bold/large title\n  // look.\n  anchorParagraph = paragraphs.items[titleIndex];\n  insertLocation = \"After\";\n  needsExplicitFormatting = true;\n} else {\n  // Fallback: title paragraph not found by exact text match, use the\n  // very first paragraph of the document.\n  anchorParagraph = paragraphs.items[0];\n  insertLocation = \"After\";\n  needsExplicitFormatting = true;\n}\n\nconst authorParagraph = anchorParagraph.insertParagraph(\n  \"Author: Craig Jones\",\n  insertLocation\n);\n\nif (needsExplicitFormatting) {\n  // Make sure the new paragraph/run explicitly carries the expected\n  // formatting when it couldn't simply inherit it from a neighboring\n  // paragraph.\n  authorParagraph.styleBuiltIn = Word.BuiltInStyleName.noSpacing;\n  authorParagraph.font.name = \"Times New Roman\";\n}\n\nawait context.sync();\n", "ps1": "# Add an \"Author: Craig Jones\" paragraph right after the \"Line.java\"\n# title paragraph (and before the existing blank paragraph that\n# follows it), matching the new-paragraph formatting (NoSpacing style,\n# Times New Roman ascii/hAnsi/cs fonts) already used by that blank\n# paragraph.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n\n# Locate the title paragraph whose text is exactly \"Line.java\".\n$titleIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs($i)\n  $txt = $p.Range.Text.Trim()\n  if ($txt -eq \"Line.java\") {\n    $titleIndex = $i\n    break\n  }\n}\n\nif ($titleIndex -eq -1) {\n  # Fallback: title paragraph not found by exact text match, use the\n  # very first paragraph of the document.\n  $titleIndex = 1\n}\n\nif ($titleIndex -lt $count) {\n  # Insert right before the paragraph that currently follows the title\n  # (normally the blank separator paragraph) so the new paragraph picks\n  # up that paragraph's formatting (No Spacing / Times New Roman, no\n  # bold/large title sizing) exactly like Word does when typing a new\n  # line in front of existing, identically-styled text.\n  $nextPara = $d.Paragraphs($titleIndex + 1)\n  $nextPara.Range.InsertParagraphBefore() | Out-Null\n  $newPara = $d.Paragraphs($titleIndex + 1)\n  $newPara.Range.Text = \"Author: Craig Jones\"\n} else {\n  # Title paragraph is the very last paragraph in the document - fall\n  # back to inserting straight after it, then fix up formatting\n  # explicitly since it would otherwise inherit the bold/large title\n  # look.\n  $titlePara = $d.Paragraphs($titleIndex)\n  $newRange = $titlePara.Range.InsertParagraphAfter()\n  $newPara = $d.Paragraphs($titleIndex + 1)\n  $newPara.Range.Text = \"Author: Craig Jones\"\n  $newPara.Style = \"No Spacing\"\n  $newPara.Range.Font.Name = \"Times New Roman\"\n}\n"}
